{"js": "// Update the two-digit-by-two-digit multiplication answers in the table.\n// Each entry is [oldText, newText]; oldText values are unique in the\n// document, so a case-sensitive whole-match search-and-replace on each\n// one is sufficient and unambiguous.\nconst replacements = [\n  [\"42\u00d730=1260\", \"60\u00d721=1260\"],\n  [\"19\u00d735=665\", \"47\u00d735=1645\"],\n  [\"35\u00d758=2030\", \"91\u00d757=5187\"],\n  [\"39\u00d773=2847\", \"25\u00d746=1150\"],\n  [\"58\u00d737=2146\", \"37\u00d773=2701\"],\n  [\"90\u00d738=3420\", \"12\u00d722=264\"],\n  [\"38\u00d780=3040\", \"97\u00d792=8924\"],\n  [\"95\u00d756=5320\", \"41\u00d752=2132\"],\n  [\"74\u00d722=1628\", \"85\u00d775=6375\"],\n  [\"78\u00d794=7332\", \"85\u00d784=7140\"],\n  [\"90\u00d733=2970\", \"81\u00d722=1782\"],\n  [\"94\u00d741=3854\", \"20\u00d770=1400\"],\n  [\"72\u00d781=5832\", \"29\u00d755=1595\"],\n  [\"81\u00d714=1134\", \"21\u00d743=903\"],\n  [\"54\u00d791=4914\", \"35\u00d743=1505\"],\n  [\"31\u00d769=2139\", \"87\u00d735=3045\"],\n  [\"57\u00d722=1254\", \"39\u00d716=624\"],\n  [\"90\u00d797=8730\", \"95\u00d780=7600\"],\n  [\"55\u00d767=3685\", \"80\u00d764=5120\"],\n  [\"85\u00d729=2465\", \"70\u00d761=4270\"],\n  [\"70\u00d760=4200\", \"43\u00d714=602\"],\n  [\"52\u00d776=3952\", \"52\u00d754=2808\"],\n  [\"87\u00d725=2175\", \"21\u00d784=1764\"],\n  [\"11\u00d716=176\", \"38\u00d761=2318\"],\n  [\"19\u00d782=1558\", \"40\u00d747=1880\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-by-two-digit multiplication answers in the table.\n# Each pair is (oldText, newText); oldText values are unique in the\n# document, so Find/Replace on the whole document content is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"42\u00d730=1260\", \"60\u00d721=1260\"),\n  @(\"19\u00d735=665\", \"47\u00d735=1645\"),\n  @(\"35\u00d758=2030\", \"91\u00d757=5187\"),\n  @(\"39\u00d773=2847\", \"25\u00d746=1150\"),\n  @(\"58\u00d737=2146\", \"37\u00d773=2701\"),\n  @(\"90\u00d738=3420\", \"12\u00d722=264\"),\n  @(\"38\u00d780=3040\", \"97\u00d792=8924\"),\n  @(\"95\u00d756=5320\", \"41\u00d752=2132\"),\n  @(\"74\u00d722=1628\", \"85\u00d775=6375\"),\n  @(\"78\u00d794=7332\", \"85\u00d784=7140\"),\n  @(\"90\u00d733=2970\", \"81\u00d722=1782\"),\n  @(\"94\u00d741=3854\", \"20\u00d770=1400\"),\n  @(\"72\u00d781=5832\", \"29\u00d755=1595\"),\n  @(\"81\u00d714=1134\", \"21\u00d743=903\"),\n  @(\"54\u00d791=4914\", \"35\u00d743=1505\"),\n  @(\"31\u00d769=2139\", \"87\u00d735=3045\"),\n  @(\"57\u00d722=1254\", \"39\u00d716=624\"),\n  @(\"90\u00d797=8730\", \"95\u00d780=7600\"),\n  @(\"55\u00d767=3685\", \"80\u00d764=5120\"),\n  @(\"85\u00d729=2465\", \"70\u00d761=4270\"),\n  @(\"70\u00d760=4200\", \"43\u00d714=602\"),\n  @(\"52\u00d776=3952\", \"52\u00d754=2808\"),\n  @(\"87\u00d725=2175\", \"21\u00d784=1764\"),\n  @(\"11\u00d716=176\", \"38\u00d761=2318\"),\n  @(\"19\u00d782=1558\", \"40\u00d747=1880\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $r = $d.Content\n  $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
